# Apply the "Add files via upload" commit:
#  - Sheet (tab) name updated to reflect the new export timestamp
#  - Reference date in column G shifted from 45628 (2024-12-02) to 45629 (2024-12-03)
#    for every data row
#  - A handful of "Saldo Previsto" / "Vl. Total" values (columns E and H)
#    were corrected for specific accounts

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the new export run
$ws.Name = "IClientBalance-20241203-083151-"

# Shift the reference date for every data row (2 through 274) from
# 45628 -> 45629 (column G)
$lastRow = 274
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = 45629
}

# Corrected balance values (columns E = "Saldo Previsto", H = "Vl. Total")
$corrections = @{
    112 = 0.81
    113 = 0.19
    173 = 426.01
    251 = 0
    258 = 0
}

foreach ($row in $corrections.Keys) {
    $value = $corrections[$row]
    $ws.Cells.Item($row, 5).Value = $value
    $ws.Cells.Item($row, 8).Value = $value
}
